# Story_Data.xlsx UI pass: the author retyped several story-text cells in
# column C, swapping the old literal "\n" line-break markers for "&"
# (plus a couple of incidental wording tweaks picked up while retyping),
# and left the selection sitting on C14 afterwards.
#
# The edits are applied in the same left-to-right "as retyped" order the
# author used (C3, C6, C7, C12, C13, C5, C4) so the workbook's shared
# string table is rebuilt in the same sequence as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = '&&전 세계에 균열이 발생했다.'
$ws.Range("C6").Value = '내가 바로 그 힘을 가진&&''헌터''다.'
$ws.Range("C7").Value = '나는 협회에 소속된 헌터다.&그것도...&헌터에 관련된 계약에 허점이 많은 시절에 계약하여 노예와 다름 없는 계약...'
$ws.Range("C12").Value = '(사이렌 소리)하… 또 균열이 발생했다.&뭔 놈의 마물이 이렇게 하루가 멀다 하고 매일 나오는지 지겹다 지겨워…'
$ws.Range("C13").Value = '(전화벨소리)왜?&(중얼거리며) 아니...학교에서 필요한 준비물이 있는데..&귀찮게 연락하지 말고, 카드로 사&(뚝 끊음&'
$ws.Range("C5").Value = '하지만,&인류는 멸망하지 않았다.&&균열에서 새어 나오는 마나를 받아들여 특별한 힘을 지니게 된 사람들 덕분이었다.'
$ws.Range("C4").Value = '차원 간의 균열이 열리며&튀어나온 온갖 마물들은&그야말로 자연재해였다.'

# Leave the active cell/selection on C14, matching the saved cursor position.
$ws.Range("C14").Select()

